$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 35; every row from 35 downward shifts down by one.
$ws.Rows.Item(35).Insert()

# Row 37 (the old row 36, "Prix Geko Bike à Bernwiller") gets its source-link cell
# (column E) updated from the generic "Inscriptions via Dossardeur" text to the
# event-specific slug "bernwiller". Touch this first so the new shared string
# "bernwiller" is created before the other brand-new strings below (keeps the
# shared-string table ordering identical to the authored workbook).
$ws.Range("E37").Value = "bernwiller"

# Fill in the newly inserted row 35 with the new calendar entry.
$ws.Range("B35").Value = "L'étape Cyclo du Tour"
$ws.Range("C35").Value = "Tour Alsace & FSGT"
$ws.Range("E35").Value = "cyclotouralsace"
$ws.Range("D35").Value = "Randonnée"
$ws.Range("A35").Value = "Sam 30 Juillet"

# Restore the view state (selection + scroll anchor) recorded in the saved file.
$win = $excel.ActiveWindow
$win.ScrollRow = 21
$win.ScrollColumn = 1
[void]$ws.Range("A36").Select()
